# Generate Report for Handback
# Two files (66424e2f-... and f5f50491-...) have been handed back and are
# now "in sync with en-US". They move to the top of the Overview sheet and
# gain populated "Latest Target File" / "Latest Handback File" / "Latest
# Handback DateTime" columns on the per-locale sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: Overview
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Hyperlinks.Delete()

$ov.Range("B2").Value2 = "Handed back: in sync with en-US"
$ov.Range("C2").Value2 = "Handed back: in sync with en-US"

$ov.Range("B3").Value2 = "Handed back: in sync with en-US"
$ov.Range("C3").Value2 = "Handed back: in sync with en-US"

$ov.Range("A4").Value2 = "9676e2ac-0a8c-4f95-a4dc-608716d2934d.md"
$ov.Range("B4").Value2 = "In Translation"
$ov.Range("C4").Value2 = "In Translation"

$ov.Range("A5").Value2 = "b3b82a64-5468-4710-a2a7-13a2a451b96b.md"
$ov.Range("B5").Value2 = "Ready for handoff"
$ov.Range("C5").Value2 = "Ready for handoff"

$ov.Range("A6").Value2 = ".localization-config"
$ov.Range("B6").Value2 = "Not to be localized"
$ov.Range("C6").Value2 = "Not to be localized"

$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/fbc9492a8f5de52d597fbca2cf3cee2514ff80d8/e2e/66424e2f-a93a-4b9a-a9b5-0af2ccca2ecb.md", "", "", "66424e2f-a93a-4b9a-a9b5-0af2ccca2ecb.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/fbc9492a8f5de52d597fbca2cf3cee2514ff80d8/e2e/f5f50491-88de-4af9-b740-a4609cb6cf84.md", "", "", "f5f50491-88de-4af9-b740-a4609cb6cf84.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/2138b65e8561f07677d3fc8953b720e9b4201693/e2e/9676e2ac-0a8c-4f95-a4dc-608716d2934d.md", "", "", "9676e2ac-0a8c-4f95-a4dc-608716d2934d.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/fbc9492a8f5de52d597fbca2cf3cee2514ff80d8/e2e/b3b82a64-5468-4710-a2a7-13a2a451b96b.md", "", "", "b3b82a64-5468-4710-a2a7-13a2a451b96b.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/fbc9492a8f5de52d597fbca2cf3cee2514ff80d8/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet: zh-cn
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Hyperlinks.Delete()

# Row 2 : 66424e2f... (now handed back)
$zh.Range("B2").Value2 = "Handed back: in sync with en-US"
$zh.Range("D2").Value2 = "2016-03-10 00:16:19"
$zh.Range("G2").Value2 = "2016-03-10 00:16:19"
$zh.Range("H2").Value2 = "Include"

# Row 3 : f5f50491... (now handed back)
$zh.Range("B3").Value2 = "Handed back: in sync with en-US"
$zh.Range("C3").Value2 = "f5f50491-88de-4af9-b740-a4609cb6cf84.946de752d542b048b9f03da93b58fbd87f1bb0ac.zh-cn.xlf"
$zh.Range("D3").Value2 = "2016-03-10 00:15:21"
$zh.Range("G3").Value2 = "2016-03-10 00:16:19"
$zh.Range("H3").Value2 = "Include"

# Row 4 : 9676e2ac... (unchanged status, still "In Translation")
$zh.Range("A4").Value2 = "9676e2ac-0a8c-4f95-a4dc-608716d2934d.md"
$zh.Range("B4").Value2 = "In Translation"
$zh.Range("C4").Value2 = "9676e2ac-0a8c-4f95-a4dc-608716d2934d.8ac0eee50d151352e064275141b9b6a2c6d010f5.zh-cn.xlf"
$zh.Range("D4").Value2 = "2016-03-10 00:12:58"
$zh.Range("G4").Value2 = "0001-01-01 00:00:00"
$zh.Range("H4").Value2 = "Include"

# Row 5 : b3b82a64... (unchanged status, "Ready for handoff")
$zh.Range("A5").Value2 = "b3b82a64-5468-4710-a2a7-13a2a451b96b.md"
$zh.Range("B5").Value2 = "Ready for handoff"
$zh.Range("C5").Value2 = "b3b82a64-5468-4710-a2a7-13a2a451b96b.97e1864d4fc2f2a2b39faba29d28a8283ce0eec1.zh-cn.xlf"
$zh.Range("D5").Value2 = "2016-03-10 00:15:21"
$zh.Range("G5").Value2 = "0001-01-01 00:00:00"
$zh.Range("H5").Value2 = "Include"

# Row 6 : .localization-config
$zh.Range("A6").Value2 = ".localization-config"
$zh.Range("B6").Value2 = "Not to be localized"
$zh.Range("D6").Value2 = "0001-01-01 00:00:00"
$zh.Range("G6").Value2 = "0001-01-01 00:00:00"
$zh.Range("H6").Value2 = "Ignored"

$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/fbc9492a8f5de52d597fbca2cf3cee2514ff80d8/e2e/66424e2f-a93a-4b9a-a9b5-0af2ccca2ecb.md", "", "", "66424e2f-a93a-4b9a-a9b5-0af2ccca2ecb.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1c5bbf324e940df8203148b52bc867721f3d90b0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/66424e2f-a93a-4b9a-a9b5-0af2ccca2ecb.49ee78572f3cd140d74cc3c1e7d88c889462ed7f.zh-cn.xlf", "", "", "66424e2f-a93a-4b9a-a9b5-0af2ccca2ecb.49ee78572f3cd140d74cc3c1e7d88c889462ed7f.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/fbc9492a8f5de52d597fbca2cf3cee2514ff80d8/e2e/66424e2f-a93a-4b9a-a9b5-0af2ccca2ecb.md", "", "", "66424e2f-a93a-4b9a-a9b5-0af2ccca2ecb.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1c5bbf324e940df8203148b52bc867721f3d90b0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/66424e2f-a93a-4b9a-a9b5-0af2ccca2ecb.49ee78572f3cd140d74cc3c1e7d88c889462ed7f.zh-cn.xlf", "", "", "66424e2f-a93a-4b9a-a9b5-0af2ccca2ecb.49ee78572f3cd140d74cc3c1e7d88c889462ed7f.zh-cn.xlf") | Out-Null

$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/fbc9492a8f5de52d597fbca2cf3cee2514ff80d8/e2e/f5f50491-88de-4af9-b740-a4609cb6cf84.md", "", "", "f5f50491-88de-4af9-b740-a4609cb6cf84.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1c5bbf324e940df8203148b52bc867721f3d90b0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/f5f50491-88de-4af9-b740-a4609cb6cf84.946de752d542b048b9f03da93b58fbd87f1bb0ac.zh-cn.xlf", "", "", "f5f50491-88de-4af9-b740-a4609cb6cf84.946de752d542b048b9f03da93b58fbd87f1bb0ac.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/fbc9492a8f5de52d597fbca2cf3cee2514ff80d8/e2e/f5f50491-88de-4af9-b740-a4609cb6cf84.md", "", "", "f5f50491-88de-4af9-b740-a4609cb6cf84.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1c5bbf324e940df8203148b52bc867721f3d90b0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/f5f50491-88de-4af9-b740-a4609cb6cf84.946de752d542b048b9f03da93b58fbd87f1bb0ac.zh-cn.xlf", "", "", "f5f50491-88de-4af9-b740-a4609cb6cf84.946de752d542b048b9f03da93b58fbd87f1bb0ac.zh-cn.xlf") | Out-Null

$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/2138b65e8561f07677d3fc8953b720e9b4201693/e2e/9676e2ac-0a8c-4f95-a4dc-608716d2934d.md", "", "", "9676e2ac-0a8c-4f95-a4dc-608716d2934d.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b557812f64bde1ada4259d08af25d0d094b86dd5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/9676e2ac-0a8c-4f95-a4dc-608716d2934d.8ac0eee50d151352e064275141b9b6a2c6d010f5.zh-cn.xlf", "", "", "9676e2ac-0a8c-4f95-a4dc-608716d2934d.8ac0eee50d151352e064275141b9b6a2c6d010f5.zh-cn.xlf") | Out-Null

$zh.Hyperlinks.Add($zh.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/fbc9492a8f5de52d597fbca2cf3cee2514ff80d8/e2e/b3b82a64-5468-4710-a2a7-13a2a451b96b.md", "", "", "b3b82a64-5468-4710-a2a7-13a2a451b96b.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1c5bbf324e940df8203148b52bc867721f3d90b0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/b3b82a64-5468-4710-a2a7-13a2a451b96b.97e1864d4fc2f2a2b39faba29d28a8283ce0eec1.zh-cn.xlf", "", "", "b3b82a64-5468-4710-a2a7-13a2a451b96b.97e1864d4fc2f2a2b39faba29d28a8283ce0eec1.zh-cn.xlf") | Out-Null

$zh.Hyperlinks.Add($zh.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/fbc9492a8f5de52d597fbca2cf3cee2514ff80d8/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet: de-de
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Hyperlinks.Delete()

# Row 2 : 66424e2f... (now handed back)
$de.Range("B2").Value2 = "Handed back: in sync with en-US"
$de.Range("D2").Value2 = "2016-03-10 00:15:28"
$de.Range("G2").Value2 = "2016-03-10 00:16:33"
$de.Range("H2").Value2 = "Include"

# Row 3 : f5f50491... (now handed back)
$de.Range("B3").Value2 = "Handed back: in sync with en-US"
$de.Range("C3").Value2 = "f5f50491-88de-4af9-b740-a4609cb6cf84.946de752d542b048b9f03da93b58fbd87f1bb0ac.de-de.xlf"
$de.Range("D3").Value2 = "2016-03-10 00:15:28"
$de.Range("G3").Value2 = "2016-03-10 00:16:33"
$de.Range("H3").Value2 = "Include"

# Row 4 : 9676e2ac... (unchanged status, still "In Translation")
$de.Range("A4").Value2 = "9676e2ac-0a8c-4f95-a4dc-608716d2934d.md"
$de.Range("B4").Value2 = "In Translation"
$de.Range("C4").Value2 = "9676e2ac-0a8c-4f95-a4dc-608716d2934d.8ac0eee50d151352e064275141b9b6a2c6d010f5.de-de.xlf"
$de.Range("D4").Value2 = "2016-03-10 00:13:27"
$de.Range("G4").Value2 = "0001-01-01 00:00:00"
$de.Range("H4").Value2 = "Include"

# Row 5 : b3b82a64... (unchanged status, "Ready for handoff")
$de.Range("A5").Value2 = "b3b82a64-5468-4710-a2a7-13a2a451b96b.md"
$de.Range("B5").Value2 = "Ready for handoff"
$de.Range("C5").Value2 = "b3b82a64-5468-4710-a2a7-13a2a451b96b.97e1864d4fc2f2a2b39faba29d28a8283ce0eec1.de-de.xlf"
$de.Range("D5").Value2 = "2016-03-10 00:15:28"
$de.Range("G5").Value2 = "0001-01-01 00:00:00"
$de.Range("H5").Value2 = "Include"

# Row 6 : .localization-config
$de.Range("A6").Value2 = ".localization-config"
$de.Range("B6").Value2 = "Not to be localized"
$de.Range("D6").Value2 = "0001-01-01 00:00:00"
$de.Range("G6").Value2 = "0001-01-01 00:00:00"
$de.Range("H6").Value2 = "Ignored"

$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/fbc9492a8f5de52d597fbca2cf3cee2514ff80d8/e2e/66424e2f-a93a-4b9a-a9b5-0af2ccca2ecb.md", "", "", "66424e2f-a93a-4b9a-a9b5-0af2ccca2ecb.md") | Out-Null
$de.Hyperlinks.Add($de.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1c5bbf324e940df8203148b52bc867721f3d90b0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/66424e2f-a93a-4b9a-a9b5-0af2ccca2ecb.49ee78572f3cd140d74cc3c1e7d88c889462ed7f.de-de.xlf", "", "", "66424e2f-a93a-4b9a-a9b5-0af2ccca2ecb.49ee78572f3cd140d74cc3c1e7d88c889462ed7f.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/fbc9492a8f5de52d597fbca2cf3cee2514ff80d8/e2e/66424e2f-a93a-4b9a-a9b5-0af2ccca2ecb.md", "", "", "66424e2f-a93a-4b9a-a9b5-0af2ccca2ecb.md") | Out-Null
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1c5bbf324e940df8203148b52bc867721f3d90b0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/66424e2f-a93a-4b9a-a9b5-0af2ccca2ecb.49ee78572f3cd140d74cc3c1e7d88c889462ed7f.de-de.xlf", "", "", "66424e2f-a93a-4b9a-a9b5-0af2ccca2ecb.49ee78572f3cd140d74cc3c1e7d88c889462ed7f.de-de.xlf") | Out-Null

$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/fbc9492a8f5de52d597fbca2cf3cee2514ff80d8/e2e/f5f50491-88de-4af9-b740-a4609cb6cf84.md", "", "", "f5f50491-88de-4af9-b740-a4609cb6cf84.md") | Out-Null
$de.Hyperlinks.Add($de.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1c5bbf324e940df8203148b52bc867721f3d90b0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/f5f50491-88de-4af9-b740-a4609cb6cf84.946de752d542b048b9f03da93b58fbd87f1bb0ac.de-de.xlf", "", "", "f5f50491-88de-4af9-b740-a4609cb6cf84.946de752d542b048b9f03da93b58fbd87f1bb0ac.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/fbc9492a8f5de52d597fbca2cf3cee2514ff80d8/e2e/f5f50491-88de-4af9-b740-a4609cb6cf84.md", "", "", "f5f50491-88de-4af9-b740-a4609cb6cf84.md") | Out-Null
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1c5bbf324e940df8203148b52bc867721f3d90b0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/f5f50491-88de-4af9-b740-a4609cb6cf84.946de752d542b048b9f03da93b58fbd87f1bb0ac.de-de.xlf", "", "", "f5f50491-88de-4af9-b740-a4609cb6cf84.946de752d542b048b9f03da93b58fbd87f1bb0ac.de-de.xlf") | Out-Null

$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/2138b65e8561f07677d3fc8953b720e9b4201693/e2e/9676e2ac-0a8c-4f95-a4dc-608716d2934d.md", "", "", "9676e2ac-0a8c-4f95-a4dc-608716d2934d.md") | Out-Null
$de.Hyperlinks.Add($de.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b557812f64bde1ada4259d08af25d0d094b86dd5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/9676e2ac-0a8c-4f95-a4dc-608716d2934d.8ac0eee50d151352e064275141b9b6a2c6d010f5.de-de.xlf", "", "", "9676e2ac-0a8c-4f95-a4dc-608716d2934d.8ac0eee50d151352e064275141b9b6a2c6d010f5.de-de.xlf") | Out-Null

$de.Hyperlinks.Add($de.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/fbc9492a8f5de52d597fbca2cf3cee2514ff80d8/e2e/b3b82a64-5468-4710-a2a7-13a2a451b96b.md", "", "", "b3b82a64-5468-4710-a2a7-13a2a451b96b.md") | Out-Null
$de.Hyperlinks.Add($de.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1c5bbf324e940df8203148b52bc867721f3d90b0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/b3b82a64-5468-4710-a2a7-13a2a451b96b.97e1864d4fc2f2a2b39faba29d28a8283ce0eec1.de-de.xlf", "", "", "b3b82a64-5468-4710-a2a7-13a2a451b96b.97e1864d4fc2f2a2b39faba29d28a8283ce0eec1.de-de.xlf") | Out-Null

$de.Hyperlinks.Add($de.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/fbc9492a8f5de52d597fbca2cf3cee2514ff80d8/.localization-config", "", "", ".localization-config") | Out-Null

Write-Host "Handback report generated."
